# Relatorio_Jogos_0x0.xlsx — refresh the "Jogos 0x0" monitoring sheet:
# new scores/teams for rows 2-6, drop the two stale rows (7-8), and shrink
# the autofilter / filter-database range to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Campeonato Islandês ---
$ws.Range("A2").Value = "Campeonato Islandês"
$ws.Range("B2").Value = "Valur Reykjavik"
$ws.Range("D2").Value = "IF Vestri"
$ws.Range("F2").Value = "40'"
$ws.Range("G2").Value = 40
$ws.Range("J2").Value = "11:42:30"

# --- Row 3: Liga Nacional de Futebol ---
$ws.Range("A3").Value = "Liga Nacional de Futebol"
$ws.Range("B3").Value = "Torpedo Moscovo"
$ws.Range("D3").Value = "SK Rotor Volgograd"
$ws.Range("F3").Value = "40'"
$ws.Range("G3").Value = 40
$ws.Range("J3").Value = "11:42:31"

# --- Row 4: Liga Principal ---
$ws.Range("A4").Value = "Liga Principal"
$ws.Range("B4").Value = "FC Gomel"
$ws.Range("D4").Value = "FC Slavia Mozyr"
$ws.Range("F4").Value = "40'"
$ws.Range("G4").Value = 40
$ws.Range("J4").Value = "11:42:32"

# --- Row 5: Primeira Liga de Futebol Profissional ---
$ws.Range("A5").Value = "Primeira Liga de Futebol Profissional"
$ws.Range("B5").Value = "Beroe Stara Zagora"
$ws.Range("D5").Value = "Levski Sófia"
$ws.Range("F5").Value = "24'"
$ws.Range("G5").Value = 24
$ws.Range("J5").Value = "11:42:33"

# --- Row 6: Torneo Federal A ---
$ws.Range("A6").Value = "Torneo Federal A"
$ws.Range("B6").Value = "CD Santamarina Tandil"
$ws.Range("D6").Value = "Deportivo Rincón"
$ws.Range("F6").Value = "36'"
$ws.Range("G6").Value = 36
$ws.Range("J6").Value = "11:42:34"

# Drop the two rows that fell out of the monitoring window.
$ws.Rows("7:8").Delete()

# Re-point the autofilter at the now-shorter table (A1:J6) without leaving
# any column filter criteria behind.
$ws.AutoFilterMode = $false
$ws.Range("A1:J6").AutoFilter(1, "Em Andamento")
$ws.Range("A1:J6").AutoFilter(1)

# Shrink the hidden _FilterDatabase defined name to match.
$fd = $wb.Names.Item(1)
$fd.RefersTo = "='Jogos 0x0'!`$A`$1:`$J`$6"
